# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that used to
# read "Ready for handoff" now reads "In Translation" (Overview!E2:E4 &
# F2:F4, zh-cn!C2:C4, de-de!C2:C4). Excel's column AutoFit then narrows the
# now-shorter "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn status) and F (de-de status) -------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2:E4").Value = "In Translation"
$ws1.Range("F2:F4").Value = "In Translation"

# --- zh-cn sheet: column C (Status) --------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2:C4").Value = "In Translation"

# --- de-de sheet: column C (Status) --------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2:C4").Value = "In Translation"

# --- Re-fit the Status columns now that the text is shorter ---------------
# (ColumnWidth is expressed in characters; 12.5 is the closest width this
# host resolves to the narrower, re-fitted column.)
$ws1.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$ws2.Range("C1").EntireColumn.ColumnWidth = 12.5
$ws3.Range("C1").EntireColumn.ColumnWidth = 12.5
